# Convert the relative redirect paths in column A/B of row 2 into full
# URLs and turn those two cells into live hyperlinks pointing at the
# same address (mirrors "Test the protocol redirect adjustments").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "https://example.com/from"
$ws.Range("B2").Value = "https://example.com/to"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://example.com/from")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://example.com/to")

# Leave the cursor where the author left it when they saved.
$null = $ws.Range("A11").Select()
